# Update the "Förändrad" (Changed) date column (C) for rows 2-33 from
# 45205 (2023-10-06) to 45206 (2023-10-07), reflecting an automatic
# update timestamp bump, as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 3).Value = 45206
}
